$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": add column X (07-jul) ---
$wsSpot = $wb.Worksheets.Item("Prix Spot")

$wsSpot.Range("W1").Copy()
$wsSpot.Range("X1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsSpot.Range("X1").Value = "07-jul"

$spotValues = @(46.35, 35.82, 38.15, 25.44, 29.2, 35.85, 40.46, 38.13, 58.71, 57.5, 37.5, 18.01, 26.99, 8.13, 8.7, 5, 3.6, 4.31, 34.47, 54.37, 50.2, 54.66, 92.27, 78.98)

for ($i = 0; $i -lt $spotValues.Length; $i++) {
    $row = $i + 2
    $wsSpot.Cells.Item($row, 24).Value = $spotValues[$i]
}

# --- Sheet "Gaz": add rows 21-22 ---
$wsGaz = $wb.Worksheets.Item("Gaz")

$wsGaz.Range("A21").NumberFormat = "@"
$wsGaz.Range("A21").Value = "2025-07-05"
$wsGaz.Range("A21").Style = "Normal"
$wsGaz.Range("B21").Value = 32.5

$wsGaz.Range("A22").NumberFormat = "@"
$wsGaz.Range("A22").Value = "2025-07-06"
$wsGaz.Range("A22").Style = "Normal"
$wsGaz.Range("B22").Value = 32.5

# --- Sheet "CO2": add rows 21-22 ---
$wsCO2 = $wb.Worksheets.Item("CO2")

$wsCO2.Range("A21").NumberFormat = "@"
$wsCO2.Range("A21").Value = "2025-07-05"
$wsCO2.Range("A21").Style = "Normal"
$wsCO2.Range("B21").Value = 70.92

$wsCO2.Range("A22").NumberFormat = "@"
$wsCO2.Range("A22").Value = "2025-07-06"
$wsCO2.Range("A22").Style = "Normal"
$wsCO2.Range("B22").Value = 70.92
